$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 15

$ws.Cells.Item($row, 1).Value = 42619.890370370369
$ws.Cells.Item($row, 2).Value = 6
$ws.Cells.Item($row, 3).Value = 58
$ws.Cells.Item($row, 4).Value = 38
$ws.Cells.Item($row, 5).Value = 58
$ws.Cells.Item($row, 6).Value = 35
$ws.Cells.Item($row, 7).Value = 16835
$ws.Cells.Item($row, 8).Value = 19481
$ws.Cells.Item($row, 9).Value = 3272
$ws.Cells.Item($row, 10).Value = 487
$ws.Cells.Item($row, 11).Value = 318
$ws.Cells.Item($row, 12).Value = 67
$ws.Cells.Item($row, 13).Value = 37
$ws.Cells.Item($row, 14).Value = "Noun"

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
